# project update version 63493
# Freezes the previously-volatile RANDBETWEEN() sample data on the "Data"
# sheet into plain literal values (rows 12-14 and 20-24), and restores the
# last-used cursor position on the "Data" sheet without changing which
# sheet tab is active (the "Charts" sheet stays the active tab).

$wb = $excel.ActiveWorkbook
$wsData = $wb.Worksheets.Item("Data")
$wsCharts = $wb.Worksheets.Item("Charts")

function Set-RowValues($sheet, $row, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $sheet.Cells.Item($row, 2 + $i).Value = $values[$i]
    }
}

Set-RowValues $wsData 12 @(127, 91, 148, 106, 136, 117, 99, 101, 83, 115, 119, 72)
Set-RowValues $wsData 13 @(84, 59, 99, 65, 97, 91, 69, 70, 63, 70, 74, 50)
Set-RowValues $wsData 14 @(55, 35, 50, 45, 63, 58, 41, 36, 15, 50, 31, 4)

Set-RowValues $wsData 20 @(147, 146, 206, 215, 262, 215, 193, 159, 194, 180, 154, 148)
Set-RowValues $wsData 21 @(204, 211, 314, 283, 358, 262, 295, 232, 226, 264, 180, 208)
Set-RowValues $wsData 22 @(116, 121, 162, 194, 214, 192, 149, 131, 173, 148, 122, 118)
Set-RowValues $wsData 23 @(146, 206, 215, 262, 215, 193, 159, 194, 180, 154, 148, 165)
Set-RowValues $wsData 24 @(142, 125, 141, 134, 136, 124, 124, 135, 131, 145, 125, 127)

# Restore the saved selection on the Data sheet (M16), then re-activate the
# Charts sheet so it remains the workbook's active tab.
$wsData.Range("M16").Select()
$wsCharts.Activate()
